$d = $word.ActiveDocument

# Find the paragraph that ends the "JQuery" list item so we can insert the
# two new technology bullets ("HTML5" and "CSS3") right after it, before the
# "Tambien se utilizara..." paragraph.
$jquery = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "JQuery") {
        $jquery = $p
    }
}

if ($jquery -eq $null) {
    throw "Could not find the 'JQuery' paragraph"
}

# Insert a new paragraph right after JQuery; it inherits JQuery's paragraph
# formatting (pStyle "Prrafodelista" + numPr list numbering + es-CO lang),
# matching the rest of the bulleted technology list.
$jquery.Range.InsertParagraphAfter()
$html5 = $jquery.Next()
$html5.Range.Text = "HTML5"

# Insert another new paragraph after the HTML5 one for CSS3, same formatting.
$html5.Range.InsertParagraphAfter()
$css3 = $html5.Next()
$css3.Range.Text = "CSS3"
